# Problem 10 draft update: fill in a couple of missing DP cells in the
# existing "ab*ac*a" table (rows 45-51) and append three new worked
# examples (pattern ".*c", pattern "a", pattern "bbbba") as new DP tables
# starting at row 54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: paste only the formatting (style) of $srcAddr onto $dstAddr
# without touching its value/content.
# ---------------------------------------------------------------------
function Copy-Style([string]$srcAddr, [string]$dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 1. Fill in two missing results in the existing "ab*ac*a" table.
# ---------------------------------------------------------------------
$ws.Range("C48").Value = 1
Copy-Style "C4" "C48"

$ws.Range("D50").Value = 1
Copy-Style "C4" "D50"

# ---------------------------------------------------------------------
# 2. New example table at row 54: pattern ".*c" vs string "ab"
# ---------------------------------------------------------------------
$ws.Range("B54:AA54").Merge()
Copy-Style "B1:AA1" "B54:AA54"
$ws.Range("B54").Value = ".*c"

$ws.Range("B55").Value = "a"
$ws.Range("C55").Value = "b"

$ws.Range("A56").Value = ".*"
$ws.Range("B56").Value = 1
$ws.Range("C56").Value = 2
Copy-Style "F39" "B56"
Copy-Style "F39" "C56"

$ws.Range("A57").Value = "c"
$ws.Range("C57").Value = 1
Copy-Style "F39" "C57"

# ---------------------------------------------------------------------
# 3. New example table at row 60: pattern "a" vs string "a"
# ---------------------------------------------------------------------
$ws.Range("B60:AA60").Merge()
Copy-Style "B1:AA1" "B60:AA60"
$ws.Range("B60").Value = "a"

$ws.Range("B61").Value = "a"

$ws.Range("A62").Value = "a"
$ws.Range("B62").Value = 1
Copy-Style "C4" "B62"

$ws.Range("A63").Value = "b*"

# ---------------------------------------------------------------------
# 4. New example table at row 66: pattern "bbbba" vs string "bbba"
# ---------------------------------------------------------------------
$ws.Range("B66:AA66").Merge()
Copy-Style "B1:AA1" "B66:AA66"
$ws.Range("B66").Value = "bbbba"

$ws.Range("B67").Value = "b"
$ws.Range("C67").Value = "b"
$ws.Range("D67").Value = "b"
$ws.Range("E67").Value = "a"

$ws.Range("A68").Value = ".*"
Copy-Style "C4" "B68"
Copy-Style "C4" "C68"
Copy-Style "C4" "D68"

$ws.Range("A69").Value = "a*"

$ws.Range("A70").Value = "a"
Copy-Style "C4" "E70"

# ---------------------------------------------------------------------
# 5. Selection / view state (matches the commit's last-touched cell).
# ---------------------------------------------------------------------
$ws.Range("E68").Select()
